$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "submission_property_type.id" row (row 43), which is an
# internal primary key without a meaning. Deleting the entire row shifts
# subsequent rows up automatically, and the sheet's used-range dimension
# shrinks from A1:D47 to A1:D46.
$ws.Rows.Item(43).Delete()

# Reflect the reviewer's resulting view state: zoomed in a bit and with
# the row that used to hold submission_type (now the last data row)
# selected.
$excel.ActiveWindow.Zoom = 130
$ws.Range("A43:XFD43").Select()
